# Bump version to 1.0.0 for first release
# Update BOM entry for SW3 (SPDT power switch) to the new part: PCM12SMTR / C&K Components

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "PCM12SMTR"
$ws.Range("F16").Value = "C&K Components"
$ws.Range("G16").Value = "OpenBCI_Wifi_Shield:PCM12SMTR"

# Move the active selection to G17, matching the saved view state after the edit
$ws.Range("G17").Select()
